$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing data but keep the formatting (bold/border/center-top style) that
# already lives on the header row (B1:J1) and the index column (A2:A11) ---
$ws.Cells.ClearContents()
$ws.Range("A11").Clear()

# --- Bring the two new columns (K, L) up to the same header style used by B1:J1 ---
$ws.Range("B1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# --- Header row (written left to right so new shared strings land in the same order
# a plain column-by-column rewrite of the sheet would produce) ---
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "Gender"
$ws.Range("G1").Value = "Birth-day"
$ws.Range("H1").Value = "Birth-month"
$ws.Range("I1").Value = "Birth-year"
$ws.Range("J1").Value = "Height"
$ws.Range("K1").Value = "Weight"
$ws.Range("L1").Value = "Calorie goal"

# --- Column B: First Name ---
$ws.Range("B2").Value = "John"
$ws.Range("B3").Value = "Steve"
$ws.Range("B4").Value = "Stan"
$ws.Range("B5").Value = "Jasmine"
$ws.Range("B6").Value = "Clark"
$ws.Range("B7").Value = "Elicina"
$ws.Range("B8").Value = "Lloyd"
$ws.Range("B9").Value = "Bob"
$ws.Range("B10").Value = "Roy"

# --- Column C: Last Name ---
$ws.Range("C2").Value = "Legend"
$ws.Range("C3").Value = "Smith"
$ws.Range("C4").Value = "Smith"
$ws.Range("C5").Value = "Mai"
$ws.Range("C6").Value = "Kent"
$ws.Range("C7").Value = "Crimea"
$ws.Range("C8").Value = "Reed"
$ws.Range("C9").Value = "Smith"
$ws.Range("C10").Value = "Harper"

# --- Column D: Email ---
$ws.Range("D2").Value = "john.legend@gmail.com"
$ws.Range("D3").Value = "steve.smith@gmail.com"
$ws.Range("D4").Value = "stan.smith@gmail.com"
$ws.Range("D5").Value = "jasmine@gmail.com"
$ws.Range("D6").Value = "superman@gmail.com"
$ws.Range("D7").Value = "elincia.crimea@gmail.com"
$ws.Range("D8").Value = "lloyd.reed@gmail.com"
$ws.Range("D9").Value = "bob.smith.1990@gmail.com"
$ws.Range("D10").Value = "roy.harper@gmail.com"

# --- Column E: Password ---
$ws.Range("E2").Value = "Dog1"
$ws.Range("E3").Value = "Stevetheboy2"
$ws.Range("E4").Value = "Stantheman1"
$ws.Range("E5").Value = "Cat2"
$ws.Range("E6").Value = "iamsuperman"
$ws.Range("E7").Value = "blesscrimea"
$ws.Range("E8").Value = "whitewolf"
$ws.Range("E9").Value = "bobsmith1990"
$ws.Range("E10").Value = "iamredarrow"

# --- Column F: Gender ---
$ws.Range("F2").Value = "Male"
$ws.Range("F3").Value = "Male"
$ws.Range("F4").Value = "Male"
$ws.Range("F5").Value = "Female"
$ws.Range("F6").Value = "Male"
$ws.Range("F7").Value = "Female"
$ws.Range("F8").Value = "Male"
$ws.Range("F9").Value = "Male"
$ws.Range("F10").Value = "Male"

# --- Column A: row index (numbers, keeps its pre-existing style) ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8

# --- Columns G:L for the numeric rows (Birth-day, Birth-month, Birth-year, Height, Weight, Calorie goal) ---
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1991
$ws.Range("J2").Value = 72
$ws.Range("K2").Value = 210
$ws.Range("L2").Value = 100

$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 1992
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 150
$ws.Range("L3").Value = 100

$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 1993
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 100

$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 1994
$ws.Range("J5").Value = 68
$ws.Range("K5").Value = 140
$ws.Range("L5").Value = 100

$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1990
$ws.Range("J6").Value = 76
$ws.Range("K6").Value = 200
$ws.Range("L6").Value = 300

$ws.Range("G7").Value = 17
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 1998
$ws.Range("J7").Value = 66
$ws.Range("K7").Value = 130
$ws.Range("L7").Value = 500

$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1997
$ws.Range("J8").Value = 71
$ws.Range("K8").Value = 160
$ws.Range("L8").Value = 600

$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 100

# --- Row 10 (Roy Harper) stores its numeric-looking fields as TEXT, not numbers ---
$ws.Range("G10:L10").NumberFormat = "@"
$ws.Range("G10").Value = "1"
$ws.Range("H10").Value = "1"
$ws.Range("I10").Value = "2000"
$ws.Range("J10").Value = "74"
$ws.Range("K10").Value = "160"
$ws.Range("L10").Value = "500"
$ws.Range("G10:L10").Style = "Normal"

Write-Output "edit applied"
